$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library preparation protocol")

# Insert a new column before column P (16th column), pushing existing
# P..AU content right to Q..AV, carrying formatting along.
$ws.Columns.Item(16).Insert()

# Size the freshly inserted column to match the new "nucleic_acid_source"
# header/value column.
$ws.Columns.Item(16).ColumnWidth = 14.5

# Row 2 (machine-readable field-name row): new column header key.
$ws.Cells.Item(2, 16).Value = "library_preparation_protocol.nucleic_acid_source"

# Row 4 (example row): example value, left in the default/general style.
$ws.Cells.Item(4, 16).Value = "library_preparation_protocol.nucleic_acid_source"
$ws.Cells.Item(4, 16).Style = "Normal"

# Row 5 (blank data-entry separator row): keep blank but match the plain
# bold/grey style used elsewhere in that row group (not the bordered one).
$ws.Cells.Item(5, 16).Style = $ws.Cells.Item(5, 1).Style

# Row 6 (human-readable description row): the actual description text,
# restoring the wrap-text style used by its neighbouring description cells.
$ws.Cells.Item(6, 16).Value = "single cell"
$ws.Cells.Item(6, 16).Style = $ws.Cells.Item(6, 13).Style
